$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Neo4j/Cypher query text that now documents the "dbExcel" query cell (A2).
# A2 already carries the wrap-text style (s="1"), so simply writing the long
# string there reproduces the wrapped, taller row seen in the target file.
$query = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Adenocarcinoma of the small intestine''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

$ws.Range("A2").Value = $query

# The long wrapped text makes Excel grow row 2 to fit its content.
$ws.Rows.Item(2).RowHeight = 101.5

# Selection moved from C7 to the B2:B5 query-column block.
$ws.Range("B2:B5").Select()
